$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so numeric-looking
# strings like "330.66" or "0.00000000328" are not reinterpreted as
# floating point numbers and re-serialized with precision artifacts.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "29.551.42"
$ws.Range("E2").Value = "  -3.01%  "

$ws.Range("D3").Value = "2.004.07"
$ws.Range("E3").Value = "  -5.31%  "

$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "330.66"
$ws.Range("E5").Value = "  -4.55%  "

$ws.Range("D6").Value = "1.006"
$ws.Range("E6").Value = "  -0.07%  "

$ws.Range("D7").Value = "0.5023"
$ws.Range("E7").Value = "  -4.09%  "

$ws.Range("D8").Value = "0.4267"
$ws.Range("E8").Value = "  -4.37%  "

$ws.Range("D9").Value = "53.81"
$ws.Range("E9").Value = "  -0.68%  "

$ws.Range("D10").Value = "0.09200"
$ws.Range("E10").Value = "  -2.27%  "

$ws.Range("D11").Value = "1.127"
$ws.Range("E11").Value = "  -4.51%  "

$ws.Range("D12").Value = "23.54"
$ws.Range("E12").Value = "  -6.90%  "

$ws.Range("D13").Value = "8.136"
$ws.Range("E13").Value = "  -6.70%  "

$ws.Range("E14").Value = "  -5.81%  "

$ws.Range("D15").Value = "1.903.54"
$ws.Range("E15").Value = "  -10.00%  "

$ws.Range("D16").Value = "96.01"
$ws.Range("E16").Value = "  -5.88%  "

$ws.Range("E17").Value = "  +0.01%  "

$ws.Range("D18").Value = "0.00001124"
$ws.Range("E18").Value = "  -3.69%  "

$ws.Range("D19").Value = "0.06653"
$ws.Range("E19").Value = "  -1.21%  "

$ws.Range("D20").Value = "19.89"
$ws.Range("E20").Value = "  -7.33%  "

$ws.Range("D21").Value = "1.005"
$ws.Range("E21").Value = "  -0.05%  "

$ws.Range("D22").Value = "5.991"
$ws.Range("E22").Value = "  -5.50%  "

$ws.Range("D23").Value = "29.544.12"
$ws.Range("E23").Value = "  -3.22%  "

$ws.Range("D24").Value = "11.99"
$ws.Range("E24").Value = "  -5.65%  "

$ws.Range("D25").Value = "2.274"
$ws.Range("E25").Value = "  -2.27%  "

$ws.Range("D26").Value = "158.99"
$ws.Range("E26").Value = "  -2.13%  "

$ws.Range("D27").Value = "20.78"
$ws.Range("E27").Value = "  -6.32%  "

$ws.Range("D28").Value = "6.556"
$ws.Range("E28").Value = "  -5.02%  "

$ws.Range("E29").Value = "  -8.63%  "

$ws.Range("D30").Value = "128.57"
$ws.Range("E30").Value = "  -4.06%  "

$ws.Range("D31").Value = "1.056"
$ws.Range("E31").Value = "  -9.11%  "

$ws.Range("D32").Value = "1.593"
$ws.Range("E32").Value = "  -10.39%  "

$ws.Range("D33").Value = "0.09965"
$ws.Range("E33").Value = "  -6.05%  "

$ws.Range("D34").Value = "5.868"
$ws.Range("E34").Value = "  -6.80%  "

$ws.Range("D35").Value = "3.776"

$ws.Range("D36").Value = "9.647"
$ws.Range("E36").Value = "  -9.23%  "

$ws.Range("D37").Value = "0.02474"
$ws.Range("E37").Value = "  -6.81%  "

$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "1.313"
$ws.Range("E38").Value = "  -2.01%  "

$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "0.06391"
$ws.Range("E39").Value = "  -7.10%  "

$ws.Range("D40").Value = "0.6575"
$ws.Range("E40").Value = "  -7.66%  "

$ws.Range("D41").Value = "11.80"
$ws.Range("E41").Value = "  -6.87%  "

$ws.Range("D42").Value = "0.2077"
$ws.Range("E42").Value = "  -7.67%  "

$ws.Range("D43").Value = "1.004"
$ws.Range("E43").Value = "  -0.12%  "

$ws.Range("D44").Value = "0.6358"
$ws.Range("E44").Value = "  -7.82%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "13.57"
$ws.Range("E45").Value = "  -7.00%  "

$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "2.217"
$ws.Range("E46").Value = "  -7.22%  "

$ws.Range("D47").Value = "1.298"
$ws.Range("E47").Value = "  -5.67%  "

$ws.Range("D48").Value = "3.525"
$ws.Range("E48").Value = "  -3.50%  "

$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.00000000328"
$ws.Range("E49").Value = "  -4.29%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.07030"
$ws.Range("E50").Value = "  -2.89%  "

$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "1.141"
$ws.Range("E51").Value = "  -5.64%  "

# Restore the default cell style on column D now that the literal
# text values are in place (keeps the workbook styling identical to
# the original, i.e. no explicit style index on the data cells).
$priceRange.Style = "Normal"
